# "upload resource for russia block" -- add a new Minigame row (id 17000006)
# for the Russia-block game: Name / WindowId / IconPath, no BgImage.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The data lives inside an Excel Table ("表1"); growing it via ListRows.Add
# extends the table ref / autoFilter (A3:E8 -> A3:E9) the same way Excel does.
$tbl = $ws.ListObjects.Item(1)
$tbl.ListRows.Add() | Out-Null

$ws.Range("A9").Value = 17000006
$ws.Range("B9").Value = "俄罗斯块"
$ws.Range("C9").Value = 1105
$ws.Range("D9").Value = "GameButton6"
# E9 (BgImage) intentionally left blank -- no background image for this entry.

# Matches the author's final cursor position in the saved workbook.
$ws.Range("D6").Select()
